$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "ShipmentTracking" numbers for rows 2-25 (column P).
# NumberFormat "@" forces these numeric-looking strings to be stored as
# text (matching the workbook's existing shared-string cells), and
# resetting Style back to "Normal" afterwards keeps the cell style index
# unchanged (these cells carry no explicit style in the source file).
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "320018179991"
$ws.Range("P2").Style = "Normal"
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "320018180002"
$ws.Range("P3").Style = "Normal"
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "320018180035"
$ws.Range("P4").Style = "Normal"
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "320018180057"
$ws.Range("P5").Style = "Normal"
$ws.Range("P6").NumberFormat = "@"
$ws.Range("P6").Value = "320018180090"
$ws.Range("P6").Style = "Normal"
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "320018180127"
$ws.Range("P7").Style = "Normal"
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "320018180150"
$ws.Range("P8").Style = "Normal"
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "320018180171"
$ws.Range("P9").Style = "Normal"
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "320018180208"
$ws.Range("P10").Style = "Normal"
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "320018180220"
$ws.Range("P11").Style = "Normal"
$ws.Range("P12").NumberFormat = "@"
$ws.Range("P12").Value = "320018180263"
$ws.Range("P12").Style = "Normal"
$ws.Range("P13").NumberFormat = "@"
$ws.Range("P13").Value = "320018180285"
$ws.Range("P13").Style = "Normal"
$ws.Range("P14").NumberFormat = "@"
$ws.Range("P14").Value = "320018180311"
$ws.Range("P14").Style = "Normal"
$ws.Range("P15").NumberFormat = "@"
$ws.Range("P15").Value = "320018180333"
$ws.Range("P15").Style = "Normal"
$ws.Range("P16").NumberFormat = "@"
$ws.Range("P16").Value = "320018180366"
$ws.Range("P16").Style = "Normal"
$ws.Range("P17").NumberFormat = "@"
$ws.Range("P17").Value = "320018180388"
$ws.Range("P17").Style = "Normal"
$ws.Range("P18").NumberFormat = "@"
$ws.Range("P18").Value = "320018180425"
$ws.Range("P18").Style = "Normal"
$ws.Range("P19").NumberFormat = "@"
$ws.Range("P19").Value = "320018180447"
$ws.Range("P19").Style = "Normal"
$ws.Range("P20").NumberFormat = "@"
$ws.Range("P20").Value = "320018180480"
$ws.Range("P20").Style = "Normal"
$ws.Range("P21").NumberFormat = "@"
$ws.Range("P21").Value = "320018180506"
$ws.Range("P21").Style = "Normal"
$ws.Range("P22").NumberFormat = "@"
$ws.Range("P22").Value = "320018180539"
$ws.Range("P22").Style = "Normal"
$ws.Range("P23").NumberFormat = "@"
$ws.Range("P23").Value = "320018180540"
$ws.Range("P23").Style = "Normal"
$ws.Range("P24").NumberFormat = "@"
$ws.Range("P24").Value = "320018180550"
$ws.Range("P24").Style = "Normal"
$ws.Range("P25").NumberFormat = "@"
$ws.Range("P25").Value = "320018180561"
$ws.Range("P25").Style = "Normal"

$ws.Range("Q22").NumberFormat = "@"
$ws.Range("Q22").Value = "$202.67"
$ws.Range("Q22").Style = "Normal"

$ws.Range("Q24").NumberFormat = "@"
$ws.Range("Q24").Value = "$248.51"
$ws.Range("Q24").Style = "Normal"
